# The uploaded workbook had a trailing/blank " " (single space) text value
# in cells F2:H2 (IDENTIFICATION_ID, IDENTIFICATION_TYPE_DESC, TAX_ID_TYPE_CODE)
# on the single data row. Clear those stray space values, leaving the cells
# truly empty while keeping their existing formatting/style intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2:H2").ClearContents()
